# Fix typo 'Cyro-EM' to 'Cryo-EM' in the storage_medium lookup sheet,
# and re-insert it in the position the corrected value moved to,
# shifting the intervening rows down by one (matching the authoritative diff).
# Also bumps the pav:createdOn timestamp recorded on the .metadata sheet.

$wb = $excel.ActiveWorkbook

$storageMedium = $wb.Worksheets.Item("storage_medium")

# Capture the original label/URI pairs for the affected rows (12-14) before
# overwriting anything, since row 12 DMSO (serum) and row 13 RNAlater both
# need to shift down by one row to make room for the corrected Cryo-EM entry.
$dmsoLabel = $storageMedium.Range("A12").Value2
$dmsoUri   = $storageMedium.Range("B12").Value2
$rnaLabel  = $storageMedium.Range("A13").Value2
$rnaUri    = $storageMedium.Range("B13").Value2

# Row 12 becomes the corrected "Cryo-EM" entry (previously on row 14).
$storageMedium.Range("A12").Value = "Cryo-EM"
$storageMedium.Range("B12").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000333"

# Row 13 becomes the old row 12 (DMSO (serum)).
$storageMedium.Range("A13").Value = $dmsoLabel
$storageMedium.Range("B13").Value = $dmsoUri

# Row 14 becomes the old row 13 (RNAlater).
$storageMedium.Range("A14").Value = $rnaLabel
$storageMedium.Range("B14").Value = $rnaUri

# Update the recorded createdOn timestamp on the .metadata sheet.
$metadata = $wb.Worksheets.Item(".metadata")
$metadata.Range("C2").Value = "2024-03-14T10:54:38-04:00"
